$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 176987
$ws.Range("C4").Value = 166944
$ws.Range("C7").Value = 5.67
$ws.Range("C8").Value = 64.86
